# Analyse des performances apres 7eme correction
# For each sheet, copy the "ETAPE 6" column (I) values/formatting into the
# "ETAPE 7" column (J), then adjust the few measurements that actually
# changed for this round of corrections.

$wb = $excel.ActiveWorkbook

# --- Sheet "LightHouse - Portable" ---
$ws1 = $wb.Worksheets.Item("LightHouse - Portable")
$ws1.Range("I4:I7").Copy($ws1.Range("J4:J7"))
$ws1.Application.CutCopyMode = $false
$ws1.Range("I4:J7").Select() | Out-Null

# --- Sheet "LightHouse - Bureau" ---
$ws2 = $wb.Worksheets.Item("LightHouse - Bureau")
$ws2.Range("I4:I7").Copy($ws2.Range("J4:J7"))
$ws2.Application.CutCopyMode = $false
$ws2.Range("J4").Value = 89
$ws2.Range("J4").Select() | Out-Null

# --- Sheet "GTmetrix - Bureau" ---
$ws3 = $wb.Worksheets.Item("GTmetrix - Bureau")
$ws3.Range("I4:I5").Copy($ws3.Range("J4:J5"))
$ws3.Application.CutCopyMode = $false
$ws3.Range("J4").Value = 99
$ws3.Range("J5").Value = 95
$ws3.Range("J4:J5").Select() | Out-Null
